$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/percentage/URL updates - Excel keeps these as text natively
$updates = @{
    "D2" = "65.381.75"
    "E2" = "  -1.35%  "
    "D3" = "3.286.41"
    "E3" = "  -0.80%  "
    "E4" = "  -0.24%  "
    "D5" = "578.81"
    "E5" = "  +3.93%  "
    "D6" = "182.22"
    "E6" = "  -2.93%  "
    "E7" = "  -0.09%  "
    "D8" = "3.280.76"
    "E8" = "  -0.67%  "
    "D9" = "0.569"
    "E9" = "  -3.15%  "
    "E10" = "  -5.43%  "
    "D11" = "0.571"
    "E11" = "  -2.73%  "
    "D12" = "46.41"
    "E12" = "  -2.41%  "
    "E13" = "  -3.04%  "
    "D14" = "629.14"
    "E14" = "  +1.74%  "
    "D15" = "3.807.44"
    "E15" = "  -0.98%  "
    "D16" = "8.39"
    "E16" = "  -3.07%  "
    "D17" = "65.526.80"
    "E17" = "  -1.09%  "
    "E18" = "  -0.03%  "
    "B19" = "Chainlink"
    "C19" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "E19" = "  -2.51%  "
    "B20" = "WrappedEther"
    "C20" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D20" = "3.284.10"
    "E20" = "  -0.86%  "
    "D21" = "10.91"
    "E21" = "  -0.89%  "
    "D22" = "0.886"
    "E22" = "  -2.40%  "
    "D23" = "17.95"
    "E23" = "  -1.62%  "
    "D24" = "99.22"
    "E24" = "  -3.83%  "
    "D25" = "4.94"
    "E25" = "  -0.65%  "
    "D26" = "3.94"
    "E26" = "  -0.17%  "
    "D27" = "2.72"
    "E27" = "  -0.77%  "
    "E28" = "  -3.31%  "
    "D29" = "30.59"
    "E29" = "  +0.86%  "
    "D30" = "8.36"
    "E30" = "  -3.58%  "
    "D31" = "6.45"
    "E31" = "  +0.05%  "
    "D32" = "562.96"
    "E32" = "  +1.05%  "
    "E33" = "  -9.46%  "
    "D34" = "10.83"
    "E34" = "  -2.25%  "
    "D35" = "3.802.56"
    "E35" = "  -1.57%  "
    "E36" = "  -2.06%  "
    "D37" = "0.999"
    "E37" = "  -0.08%  "
    "D38" = "55.82"
    "E38" = "  -2.79%  "
    "E39" = "  -2.25%  "
    "E40" = "  +6.16%  "
    "D41" = "32.41"
    "E41" = "  -4.73%  "
    "B42" = "PEPE"
    "C42" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D42" = "0.0₃0681"
    "E42" = "  -6.42%  "
    "B43" = "Stacks"
    "C43" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D43" = "3.12"
    "E43" = "  -6.23%  "
    "D44" = "2.57"
    "E44" = "  -6.07%  "
    "E45" = "  -2.42%  "
    "D46" = "0.0403"
    "E46" = "  -4.15%  "
    "D47" = "3.02"
    "E47" = "  -6.90%  "
    "E49" = "  -2.52%  "
    "E50" = "  -3.52%  "
    "D51" = "130.39"
    "E51" = "  +6.01%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Numeric-looking price strings where Excel would normally coerce to a
# Double and lose the exact trailing/leading zero formatting (e.g. "0.330"
# -> 0.33). Force the cell to Text format first, write the literal string,
# then drop back to the Normal style so no extra formatting sticks around.
$textForced = @{
    "D13" = "0.0000263"
    "D45" = "0.330"
    "D50" = "2.50"
    "D19" = "17.60"
}

foreach ($addr in $textForced.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForced[$addr]
    $cell.Style = "Normal"
}
